# Updates cryptos list values (price + 1h volume%) to the latest scrape,
# and reorders two rows: row 50 becomes "NEARProtocol", row 51 becomes "Cronos"
# (the previous "Decentraland" row is replaced).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Plain assignment would let Excel auto-convert number-looking strings
    # (e.g. "1.000", "0.05956") into numeric values, which would both change
    # the stored type and lose the exact original-text formatting (trailing
    # zeros, etc). Forcing a text number-format keeps the literal string,
    # then resetting the style back to Normal avoids leaving a stray
    # "@"-formatted style on the cell (cells in this sheet carry no style).
    $c = $ws.Range($cell)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "30.294.03"
$ws.Range("E2").Value = "  -3.42%  "

# Row 3
$ws.Range("D3").Value = "1.929.24"
$ws.Range("E3").Value = "  -3.92%  "

# Row 4
Set-TextValue "D4" "1.000"
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
Set-TextValue "D5" "248.58"
$ws.Range("E5").Value = "  -4.06%  "

# Row 6
Set-TextValue "D6" "0.7246"
$ws.Range("E6").Value = "  -8.06%  "

# Row 7
Set-TextValue "D7" "0.9997"
$ws.Range("E7").Value = "  +0.03%  "

# Row 8
Set-TextValue "D8" "0.3323"
$ws.Range("E8").Value = "  -7.73%  "

# Row 9
Set-TextValue "D9" "28.34"
$ws.Range("E9").Value = "  -1.32%  "

# Row 10
Set-TextValue "D10" "0.06918"
$ws.Range("E10").Value = "  -2.59%  "

# Row 11
Set-TextValue "D11" "0.8034"
$ws.Range("E11").Value = "  -6.25%  "

# Row 12
Set-TextValue "D12" "0.08065"
$ws.Range("E12").Value = "  -1.07%  "

# Row 13
$ws.Range("D13").Value = "1.930.83"
$ws.Range("E13").Value = "  -3.83%  "

# Row 14
Set-TextValue "D14" "5.415"
$ws.Range("E14").Value = "  -3.76%  "

# Row 15
Set-TextValue "D15" "94.78"
$ws.Range("E15").Value = "  -6.51%  "

# Row 16
Set-TextValue "D16" "14.55"
$ws.Range("E16").Value = "  -3.36%  "

# Row 17
$ws.Range("D17").Value = "30.255.82"
$ws.Range("E17").Value = "  -3.52%  "

# Row 18
Set-TextValue "D18" "0.000008327"
$ws.Range("E18").Value = "  +4.28%  "

# Row 19
Set-TextValue "D19" "253.50"
$ws.Range("E19").Value = "  -8.29%  "

# Row 20
Set-TextValue "D20" "5.806"
$ws.Range("E20").Value = "  -2.15%  "

# Row 21
$ws.Range("D21").Value = "2.175.44"
$ws.Range("E21").Value = "  -4.02%  "

# Row 22
Set-TextValue "D22" "1.0000"
$ws.Range("E22").Value = "  +0.03%  "

# Row 23
Set-TextValue "D23" "0.9995"
$ws.Range("E23").Value = "  +0.00%  "

# Row 24
Set-TextValue "D24" "6.872"
$ws.Range("E24").Value = "  -4.65%  "

# Row 25
Set-TextValue "D25" "9.742"
$ws.Range("E25").Value = "  -3.84%  "

# Row 26
Set-TextValue "D26" "159.44"
$ws.Range("E26").Value = "  -3.27%  "

# Row 27
Set-TextValue "D27" "2.437"
$ws.Range("E27").Value = "  +1.32%  "

# Row 28
Set-TextValue "D28" "19.20"
$ws.Range("E28").Value = "  -4.20%  "

# Row 29
Set-TextValue "D29" "0.1342"
$ws.Range("E29").Value = "  -11.80%  "

# Row 30
$ws.Range("E30").Value = "  -4.85%  "

# Row 31
Set-TextValue "D31" "1.335"
$ws.Range("E31").Value = "  -1.64%  "

# Row 32
Set-TextValue "D32" "4.406"
$ws.Range("E32").Value = "  -5.37%  "

# Row 33
$ws.Range("E33").Value = "  -5.20%  "

# Row 34
Set-TextValue "D34" "0.05113"
$ws.Range("E34").Value = "  -3.06%  "

# Row 35
$ws.Range("E35").Value = "  -0.63%  "

# Row 36
Set-TextValue "D36" "0.7414"
$ws.Range("E36").Value = "  -3.96%  "

# Row 37
Set-TextValue "D37" "2.732"
$ws.Range("E37").Value = "  -2.53%  "

# Row 38
Set-TextValue "D38" "0.01980"
$ws.Range("E38").Value = "  -1.93%  "

# Row 39
Set-TextValue "D39" "2.829"
$ws.Range("E39").Value = "  -3.83%  "

# Row 40
Set-TextValue "D40" "6.618"
$ws.Range("E40").Value = "  -1.85%  "

# Row 41
Set-TextValue "D41" "79.07"
$ws.Range("E41").Value = "  -2.74%  "

# Row 42
Set-TextValue "D42" "0.4473"
$ws.Range("E42").Value = "  -6.24%  "

# Row 43
Set-TextValue "D43" "2.000"
$ws.Range("E43").Value = "  -7.83%  "

# Row 44
Set-TextValue "D44" "0.9997"
$ws.Range("E44").Value = "  +0.02%  "

# Row 45
Set-TextValue "D45" "0.8368"
$ws.Range("E45").Value = "  -2.53%  "

# Row 46
Set-TextValue "D46" "102.43"
$ws.Range("E46").Value = "  -4.12%  "

# Row 47
Set-TextValue "D47" "9.756"
$ws.Range("E47").Value = "  -2.48%  "

# Row 48
Set-TextValue "D48" "7.301"
$ws.Range("E48").Value = "  -6.22%  "

# Row 49
Set-TextValue "D49" "36.57"
$ws.Range("E49").Value = "  -1.38%  "

# Row 50
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D50" "1.485"
$ws.Range("E50").Value = "  -0.71%  "

# Row 51
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D51" "0.05956"
$ws.Range("E51").Value = "  -0.49%  "
